# Add "yes" to the selected_for_analysis (column F) for the rows that
# represent the use_case_name feature outputs, and clear it for row 7
# (which is no longer considered selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToSet = @(3,9,10,11,12,13,14,15,16,17,18,19,20,21,22,37,40,52,60,61,62,63,64,65)
foreach ($r in $rowsToSet) {
    $ws.Range("F$r").Value = "yes"
}

$rowsToClear = @(7)
foreach ($r in $rowsToClear) {
    $ws.Range("F$r").ClearContents()
}

# Update the scroll position (best effort, panes stay frozen at row 1)
# and active selection to match where the author ended up after making
# the edits.
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1

$ws.Range("F66").Select()
